# Swap the presentation's applied theme colours back to the stock "Office
# Theme" palette (the deck currently carries the "Integral" theme on its
# single slide master / theme part). The companion "Office Theme" colours
# that used to live only on the (COM-unreachable) notes-master theme part
# become the colours actually applied to the slide master, matching the
# commit's theme1.xml <-> theme2.xml colour-scheme swap.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

# COM ColorScheme.Colors(n) is 1-based and in clrScheme document order:
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#  9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB values are stored as 0xBBGGRR (standard VBA/COM RGB()).

$officeThemeColors = @(
    0x000000,  # dk1     000000
    0xFFFFFF,  # lt1     FFFFFF
    0x6A5444,  # dk2     44546A
    0xE6E6E7,  # lt2     E7E6E6
    0xD59B5B,  # accent1 5B9BD5
    0x317DED,  # accent2 ED7D31
    0xA5A5A5,  # accent3 A5A5A5
    0x00C0FF,  # accent4 FFC000
    0xC47244,  # accent5 4472C4
    0x47AD70,  # accent6 70AD47
    0xC16305,  # hlink   0563C1
    0x724F95   # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $scheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
